$wb = $excel.ActiveWorkbook

# --- Settings sheet: default values for the new tray/exe app behaviour ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B1").Value = 5
$wsSettings.Range("B4").Value = 1
$wsSettings.Range("B5").Value = 1

# --- Settings_recources sheet: image viewer pack/step setting ---
$wsRes = $wb.Worksheets.Item("Settings_recources")
$wsRes.Range("B20").Value = "40"

# --- Task_settings sheet: append new maintenance-run history entries ---
$wsTasks = $wb.Worksheets.Item("Task_settings")

$wsTasks.Range("G1").Value = "|||Datum provedení: 17.12.2024 16:12:19||Zkontrolováno: 161 souborů||Starších:      153 souborů||Smazáno:       103 souborů|||Datum provedení: 18.12.2024 14:14:10||Zkontrolováno: 108 souborů||Starších: 100 souborů||Smazáno: 50 souborů|||Datum provedení: 20.12.2024 12:00:23||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 08:49:27||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 08:49:43||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 12:00:19||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 14:21:35||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů"

$wsTasks.Range("G2").Value = "|||Datum provedení: 20.12.2024 12:00:14||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 12:00:10||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů"

$wsTasks.Range("G3").Value = "|||Datum provedení: 18.12.2024 14:13:13||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 18.12.2024 14:33:24||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 08:55:50||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 08:56:02||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 09:18:00||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 09:36:54||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 10:06:31||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 10:34:32||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů|||Datum provedení: 03.01.2025 11:02:39||Zkontrolováno: 58 souborů||Starších: 50 souborů||Smazáno: 0 souborů"
